# Append new scrape rows (2025-12-02 12:39:15 JST run) into the
# "ランサーズ" sheet (first worksheet) and refresh the timestamp on the
# rows that were already present.
#
# New rows are inserted at sheet rows 3, 6 and 9 (pushing the existing
# data down), then every data row (2-10) is rewritten in full so the
# sheet ends up with the same content/order as the target workbook.
# Hyperlinks are rebuilt from scratch afterwards because row-insert does
# not keep the <hyperlinks> collection aligned with the shifted cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove existing hyperlinks up front; they get re-created at the end
# once all rows are in their final place. ---------------------------------
$ws.Hyperlinks.Delete()

# --- Insert the three brand-new rows (done low-to-high on final index so
# each insert only has to push the rows below it). ------------------------
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(9).Insert()

# --- Column H got wider (12 -> 19 chars). Excel's ColumnWidth setter adds
# a constant ~0.8333 padding relative to the stored OOXML width, so back
# that off to land on an on-disk width of exactly 19. ---------------------
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668

$timestamp = "2025-12-02 12:39:15"

# --- Row 2 (unchanged job, just the refreshed timestamp) ------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

# --- Row 3 (NEW) ------------------------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "Reactで作成されたシステム開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445968"
$ws.Range("G3").Value = 235
$ws.Range("H3").Value = "🔥React ◆開発,システム開発"

# --- Row 4 (was row 3) ------------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "管理システムの開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5445265"
$ws.Range("G4").Value = 103
$ws.Range("H4").Value = "◆開発 ◇管理"

# --- Row 5 (was row 4) ------------------------------------------------------
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【急募】iPhone用電子黒板アプリ開発者を探しています!"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5445417"
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# --- Row 6 (NEW) ------------------------------------------------------------
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "静止画のみのデジタルサイネージ向けCMS開発の依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5445947"
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = "◆開発"

# --- Row 7 (was row 5) ------------------------------------------------------
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "賃貸保証業の管理システム構築依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5445528"
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = "◇管理"

# --- Row 8 (was row 6) ------------------------------------------------------
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【UTAGE構築代行】各種初期設定・ステップ配信・会員サイトの作成など"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5445495"
$ws.Range("G8").Value = 38
$ws.Range("H8").Value = "◇サイト"

# --- Row 9 (NEW) ------------------------------------------------------------
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【至急】WordPressの画像ギャラリー表示不具合解決依頼"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5445721"
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = "○WordPress"

# --- Row 10 (was row 7) -----------------------------------------------------
$ws.Range("A10").Value = $timestamp
$ws.Range("B10").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5445466"
$ws.Range("G10").Value = 25
# (row 10 has no H value, matching the source data)

# --- Re-create the hyperlinks for every URL cell, in row order, so they
# come back out as rId1..rId9 aligned with F2..F10. -------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445968")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445265")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5445417")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5445947")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5445528")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5445495")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5445721")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5445466")

# Hyperlinks.Add() registers a second ("...applyFont=1") copy of the
# Hyperlink cell style; reapply the named "Hyperlink" style so every link
# cell keeps using the original style slot (s="1") like the source file.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
